$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.648.35'
$ws.Range('E2').Value = '  -0.89%  '

$ws.Range('D3').Value = '2.332.51'
$ws.Range('E3').Value = '  +2.26%  '

$ws.Range('E4').Value = '  +0.05%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '232.79'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.85%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.649'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +2.14%  '

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '65.61'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +3.10%  '

$ws.Range('E8').Value = '  -0.03%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.459'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +2.57%  '

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0970'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -3.85%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '56.54'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -0.80%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '26.86'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.16%  '

$ws.Range('D13').Value = '2.684.07'
$ws.Range('E13').Value = '  +2.85%  '

$ws.Range('E14').Value = '  -1.32%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '15.44'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -1.52%  '

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '6.19'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +0.32%  '

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.849'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +0.79%  '

$ws.Range('D18').Value = '2.338.72'
$ws.Range('E18').Value = '  +2.95%  '

$ws.Range('D19').Value = '43.637.88'
$ws.Range('E19').Value = '  -0.57%  '

$ws.Range('D20').Value = '0.0₃0978'
$ws.Range('E20').Value = '  -1.99%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '74.03'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +0.37%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '6.24'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +2.46%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '248.97'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -1.75%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -0.09%  '

$ws.Range('B25').Value = 'WEMIXToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '3.80'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +13.50%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -0.77%  '

$ws.Range('E27').Value = '  -2.62%  '

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '9.92'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -1.14%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '22.29'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +6.95%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '174.96'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.72%  '

$ws.Range('E31').Value = '  +4.73%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.129'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -7.03%  '

$ws.Range('E33').Value = '  +0.70%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '5.00'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +3.96%  '

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0687'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -1.60%  '

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '4.98'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +1.72%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.46'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +5.76%  '

$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '6.54'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -0.33%  '

$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '3.62'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -5.09%  '

$ws.Range('E40').Value = '  -1.84%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '9.03'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +9.72%  '

$ws.Range('E42').Value = '  -0.11%  '

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '18.12'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +2.02%  '

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.17'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +8.41%  '

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '99.08'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +0.73%  '

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.0955'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -2.01%  '

$ws.Range('E47').Value = '  -0.14%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '4.33'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -1.21%  '

$ws.Range('D49').Value = '1.446.07'
$ws.Range('E49').Value = '  -0.14%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.32'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +0.93%  '

$ws.Range('B51').Value = 'Celestia'
$ws.Range('C51').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '9.98'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -3.40%  '
